# Insert a new weekly data row right before existing row 284 ("Rabanito" /
# Vega Central Mapocho de Santiago sheet). Excel shifts every row from 284
# down through 402 one row lower (284->285, ..., 401->402, 402->403), so we
# recreate that by inserting a blank row at 284, filling it with a copy of
# the row that lands on 285 (i.e. the old row 284), and then stamping the
# new row's date (column D) with the new reading's date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 284:402 down to 285:403, leaving a blank row 284.
$ws.Rows("284:284").Insert()

# Populate the new row 284 with the same data as the row now sitting at
# 285 (the original row 284), then overwrite just the date.
$ws.Range("A285:R285").Copy($ws.Range("A284:R284"))
$ws.Cells.Item(284, 4).Value = 45006
